$d = $word.ActiveDocument

# Namespace used for the raw XML fragments spliced in via Range.InsertXML.
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Run properties shared by every run in this paragraph (Courier New, black,
# 10.5pt, es-ES east-asian language) - identical to the original run's rPr,
# reused for every new run we splice in.
$rPr = '<w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:eastAsia="es-ES"/></w:rPr>'

# 1) Drop "considerado como el escritor " right after "soldado español ".
$find = $d.Content.Find
$find.Execute("considerado como el escritor ", $false, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# 2) Locate the remaining tail of the sentence, from "más importante" through
#    to the final period - this whole span gets replaced.
$tail = "más importante de la literatura en español tras el éxito de su obra El ingenioso Hidalgo don Quijote de la Mancha, considerada por gran número de especialistas, como la primera novela moderna y una de las más importantes de la literatura universal."
$r = $d.Range(0, 0)
$r.Find.Execute($tail, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# 3) Re-seat the (hidden, auto-managed) "_GoBack" bookmark so it sits right
#    before the replacement text instead of alone in the trailing paragraph;
#    Bookmarks.Add with the existing name moves it and clears the old spot.
$bmRange = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 4) Replace the tail span with the new runs, including the two "words" that
#    got typed in and wrapping each in spell-check proofErr markers.
$xml = '<w:p ' + $wns + '>' `
  + '<w:r>' + $rPr + '<w:t xml:space="preserve">más importante de la literatura en español tras el éxito de su obra El ingenioso Hidalgo don </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r>' + $rPr + '<w:t>adasdasdasd</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r>' + $rPr + '<w:t>adasdas</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' `
  + '<w:proofErr w:type="spellStart"/>' `
  + '<w:r>' + $rPr + '<w:t>asdad</w:t></w:r>' `
  + '<w:proofErr w:type="spellEnd"/>' `
  + '</w:p>'
$r.InsertXML($xml)

Write-Host "Final content:" $d.Content.Text
